# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
# for rows 2-51 on Sheet1, per the Sun Nov  5 21:54:55 UTC 2023 GitHub
# Actions data refresh.
#
# Price values are forced to remain plain text (NumberFormat "@" while
# assigning, then reset to the default "Normal" style) so numeric-looking
# strings such as "245.14" or "52.70" are not silently converted into
# Excel numbers, matching the source data which stores these as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.125.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.890.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.28%  "
$ws.Range("E6").Value = "  +5.72%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.347"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0712"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0993"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.168.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.693"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.892.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.198.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0815"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "239.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("E26").Value = "  +21.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0562"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.934"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("E36").Value = "  -5.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0620"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "88.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.333.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +37.22%  "
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.078.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.34%  "
